# Refresh the "cryptos" price/volume table with the latest scrape.
# Mirrors the nightly GitHub Actions job: for every tracked coin, update
# its Price (D) and Volume(1h) (E) columns; two coins (Injective Protocol /
# Trust Wallet Token) also swapped ranking position this run, so their
# Coin name (B) and Link (C) columns move too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = '36.992.32';   E = '  -1.11%  ' },
    @{ Row = 3;  D = '2.012.91';    E = '  -1.86%  ' },
    @{ Row = 4;  E = '  +0.21%  ' },
    @{ Row = 5;  D = '225.47';      E = '  -1.86%  ' },
    @{ Row = 6;  D = '0.605';       E = '  -1.48%  ' },
    @{ Row = 7;  E = '  +0.05%  ' },
    @{ Row = 8;  D = '54.43';       E = '  -5.06%  ' },
    @{ Row = 9;  D = '0.374';       E = '  -3.38%  ' },
    @{ Row = 10; D = '0.0777';      E = '  -2.67%  ' },
    @{ Row = 11; E = '  -4.98%  ' },
    @{ Row = 12; D = '2.312.53';    E = '  -1.80%  ' },
    @{ Row = 13; D = '13.98';       E = '  -5.61%  ' },
    @{ Row = 14; D = '19.91';       E = '  -4.31%  ' },
    @{ Row = 15; D = '5.19';        E = '  -2.54%  ' },
    @{ Row = 16; D = '0.736';       E = '  -3.16%  ' },
    @{ Row = 17; D = '2.013.32';    E = '  -2.51%  ' },
    @{ Row = 18; D = '36.877.59';   E = '  -1.13%  ' },
    @{ Row = 19; D = '6.35';        E = '  +4.44%  ' },
    @{ Row = 20; D = '68.27';       E = '  -2.05%  ' },
    @{ Row = 21; D = '0.0₃0810';    E = '  -2.72%  ' },
    @{ Row = 22; D = '221.61';      E = '  -2.23%  ' },
    @{ Row = 23; E = '  -0.14%  ' },
    @{ Row = 24; E = '  +2.03%  ' },
    @{ Row = 25; D = '2.16';        E = '  -6.08%  ' },
    @{ Row = 26; D = '164.97';      E = '  -2.34%  ' },
    @{ Row = 27; D = '9.04';        E = '  -6.57%  ' },
    @{ Row = 28; E = '  -2.23%  ' },
    @{ Row = 29; D = '18.50';       E = '  -2.71%  ' },
    @{ Row = 30; D = '1.29';        E = '  -6.71%  ' },
    @{ Row = 31; E = '  -2.46%  ' },
    @{ Row = 32; D = '4.44';        E = '  -2.72%  ' },
    @{ Row = 33; D = '0.0599';      E = '  -3.06%  ' },
    @{ Row = 34; D = '4.43';        E = '  -3.76%  ' },
    @{ Row = 35; D = '2.32';        E = '  -5.44%  ' },
    @{ Row = 36; D = '1.87';        E = '  +1.37%  ' },
    @{ Row = 37; E = '  +0.11%  ' },
    @{ Row = 38; E = '  -5.73%  ' },
    @{ Row = 39; D = '5.36';        E = '  +0.02%  ' },
    @{ Row = 40; D = '1.452.61';    E = '  -2.22%  ' },
    @{ Row = 41; D = '94.65';       E = '  -2.02%  ' },
    @{ Row = 42; D = '0.0211';      E = '  -5.46%  ' },
    @{ Row = 43; D = '2.79';        E = '  -3.77%  ' },
    @{ Row = 44; E = '  -3.96%  ' },
    @{ Row = 45; B = 'TrustWalletToken';  C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt';   D = '1.13';  E = '  -3.65%  ' },
    @{ Row = 46; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '16.03'; E = '  -7.19%  ' },
    @{ Row = 47; D = '7.13';        E = '  -0.51%  ' },
    @{ Row = 48; D = '0.998';       E = '  -2.68%  ' },
    @{ Row = 49; D = '2.91';        E = '  -0.46%  ' },
    @{ Row = 50; D = '2.206.10';    E = '  -1.62%  ' },
    @{ Row = 51; D = '3.54';        E = '  -10.34%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        # Price column holds plain text in this sheet (e.g. "2.013.32",
        # "0.0599"). Force text formatting first so values that parse as a
        # plain decimal (like "225.47") aren't silently coerced into a
        # Number cell / float-rounded by Excel's smart entry.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Range("E$row").Value = $u.E
    }
}
